# Ask remove hub data
# Adds a new block of localization strings (part/var/de/en) to the
# language table used to ask the user whether HUB data should be removed,
# kept, or whether the HUB plugin should just be deactivated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 255 - the confirmation question itself (taller row, wraps to 2 lines)
$ws.Range("A255").Value = "confirm"
$ws.Range("B255").Value = "remove_hub_data"
$ws.Range("C255").Value = "Möchten Sie die HUB-Daten auch entfernen? Falls Sie später einmal auf HUB2 wechseln möchten, können die Daten übernommen werden. Allenfalls möchten Sie das HUB-Plugin nur deaktivieren?"
$ws.Range("D255").Value = "Do you want to remove the HUB data as well? If you want to switch to HUB2 later, the data can be transferred. At most, you just want to disable the HUB plugin?"
$ws.Rows.Item(255).RowHeight = 30

# Row 256 - generic cancel label (reuses existing "Abbrechen" translation)
$ws.Range("B256").Value = "cancel"
$ws.Range("C256").Value = "Abbrechen"
$ws.Range("D256").Value = "Cancel"
$ws.Rows.Item(256).RowHeight = 15.75

# Row 257 - "remove HUB data" button
$ws.Range("A257").Value = "remove"
$ws.Range("B257").Value = "hub_data"
$ws.Range("C257").Value = "Entferne HUB-Daten"
$ws.Range("D257").Value = "Remove HUB data"
$ws.Rows.Item(257).RowHeight = 15.75

# Row 258 - "keep HUB data" button
$ws.Range("A258").Value = "keep"
$ws.Range("B258").Value = "hub_data"
$ws.Range("C258").Value = "HUB-Daten behalten"
$ws.Range("D258").Value = "Keep HUB data"
$ws.Rows.Item(258).RowHeight = 15.75

# Row 259 - "just deactivate HUB plugin" button
$ws.Range("A259").Value = "deactivate"
$ws.Range("B259").Value = "hub"
$ws.Range("C259").Value = "HUB-Plugin nur deaktivieren"
$ws.Range("D259").Value = "Just deactivate HUB plugin"
$ws.Rows.Item(259).RowHeight = 15.75

# Row 260 - generic "HUB data" label
$ws.Range("A260").Value = "hub"
$ws.Range("B260").Value = "data"
$ws.Range("C260").Value = "HUB-Daten"
$ws.Range("D260").Value = "HUB data"
$ws.Rows.Item(260).RowHeight = 15.75

# Row 261 - message shown after the HUB data was removed
$ws.Range("A261").Value = "msg"
$ws.Range("B261").Value = "removed_hub_data"
$ws.Range("C261").Value = "Die HUB-Daten wurden auch entfernt!"
$ws.Range("D261").Value = "The HUB data was also removed!"
$ws.Rows.Item(261).RowHeight = 15.75

# Row 262 - message shown after the HUB data was kept
$ws.Range("A262").Value = "msg"
$ws.Range("B262").Value = "kept_hub_data"
$ws.Range("C262").Value = "Die HUB-Daten wurden behalten!"
$ws.Range("D262").Value = "The HUB data was kept!"
$ws.Rows.Item(262).RowHeight = 15.75

# Move the selection to the new last row, as in the edited workbook.
$ws.Range("A263").Select()
